# This workbook tracks daily "Plátano" (banana/plantain) price observations
# for "Vega Central Mapocho de Santiago". A new reporting day (date serial
# 44585) worth 4 rows of data needs to be inserted right before the existing
# row 874, pushing all the subsequent rows down by 4 (935 -> 939 total rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at the position of the current row 874. This shifts
# the existing rows 874:935 down to 878:939, carrying over their formatting
# (including the date number format in column D) automatically.
$ws.Rows("874:877").Insert()

# Now populate the 4 freshly inserted rows (874:877) with the new
# observations for date serial 44585 (2022-01-24).

# Row 874: Barraganete / Primera
$ws.Cells.Item(874, 1).Value  = 9
$ws.Cells.Item(874, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(874, 3).Value  = "Metropolitana"
$ws.Cells.Item(874, 4).Value  = 44585
$ws.Cells.Item(874, 5).Value  = 13
$ws.Cells.Item(874, 6).Value  = "Fruta"
$ws.Cells.Item(874, 7).Value  = 100108
$ws.Cells.Item(874, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(874, 9).Value  = 100108006
$ws.Cells.Item(874, 10).Value = "Plátano"
$ws.Cells.Item(874, 11).Value = "Barraganete"
$ws.Cells.Item(874, 12).Value = "Primera"
$ws.Cells.Item(874, 13).Value = 240
$ws.Cells.Item(874, 14).Value = 23000
$ws.Cells.Item(874, 15).Value = 23000
$ws.Cells.Item(874, 16).Value = 23000
$ws.Cells.Item(874, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(874, 18).Value = "Ecuador"
$ws.Cells.Item(874, 19).Value = 1150
$ws.Cells.Item(874, 20).Value = 20

# Row 875: Sin especificar / Pintón
$ws.Cells.Item(875, 1).Value  = 9
$ws.Cells.Item(875, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(875, 3).Value  = "Metropolitana"
$ws.Cells.Item(875, 4).Value  = 44585
$ws.Cells.Item(875, 5).Value  = 13
$ws.Cells.Item(875, 6).Value  = "Fruta"
$ws.Cells.Item(875, 7).Value  = 100108
$ws.Cells.Item(875, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(875, 9).Value  = 100108006
$ws.Cells.Item(875, 10).Value = "Plátano"
$ws.Cells.Item(875, 11).Value = "Sin especificar"
$ws.Cells.Item(875, 12).Value = "Pintón"
$ws.Cells.Item(875, 13).Value = 780
$ws.Cells.Item(875, 14).Value = 10000
$ws.Cells.Item(875, 15).Value = 11000
$ws.Cells.Item(875, 16).Value = 10487
$ws.Cells.Item(875, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(875, 18).Value = "Ecuador"
$ws.Cells.Item(875, 19).Value = 524
$ws.Cells.Item(875, 20).Value = 20

# Row 876: Sin especificar / Primera Maduro
$ws.Cells.Item(876, 1).Value  = 9
$ws.Cells.Item(876, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(876, 3).Value  = "Metropolitana"
$ws.Cells.Item(876, 4).Value  = 44585
$ws.Cells.Item(876, 5).Value  = 13
$ws.Cells.Item(876, 6).Value  = "Fruta"
$ws.Cells.Item(876, 7).Value  = 100108
$ws.Cells.Item(876, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(876, 9).Value  = 100108006
$ws.Cells.Item(876, 10).Value = "Plátano"
$ws.Cells.Item(876, 11).Value = "Sin especificar"
$ws.Cells.Item(876, 12).Value = "Primera Maduro"
$ws.Cells.Item(876, 13).Value = 780
$ws.Cells.Item(876, 14).Value = 12000
$ws.Cells.Item(876, 15).Value = 13000
$ws.Cells.Item(876, 16).Value = 12487
$ws.Cells.Item(876, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(876, 18).Value = "Ecuador"
$ws.Cells.Item(876, 19).Value = 624
$ws.Cells.Item(876, 20).Value = 20

# Row 877: Sin especificar / Primera Pintón
$ws.Cells.Item(877, 1).Value  = 9
$ws.Cells.Item(877, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(877, 3).Value  = "Metropolitana"
$ws.Cells.Item(877, 4).Value  = 44585
$ws.Cells.Item(877, 5).Value  = 13
$ws.Cells.Item(877, 6).Value  = "Fruta"
$ws.Cells.Item(877, 7).Value  = 100108
$ws.Cells.Item(877, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(877, 9).Value  = 100108006
$ws.Cells.Item(877, 10).Value = "Plátano"
$ws.Cells.Item(877, 11).Value = "Sin especificar"
$ws.Cells.Item(877, 12).Value = "Primera Pintón"
$ws.Cells.Item(877, 13).Value = 840
$ws.Cells.Item(877, 14).Value = 13000
$ws.Cells.Item(877, 15).Value = 14000
$ws.Cells.Item(877, 16).Value = 13548
$ws.Cells.Item(877, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(877, 18).Value = "Ecuador"
$ws.Cells.Item(877, 19).Value = 677
$ws.Cells.Item(877, 20).Value = 20
